# ExchangeRates.xlsx edit: rename sheet, update date-range caption, insert a
# new leading rate row (30/12/2018), and append 6 trailing rate rows through
# 31/12/2019 (start of the "interest support" data extension).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the worksheet (scratch name while the new rows are wired up) ---
$ws.Name = "_tmp_exchanges"

# --- 2. Update the "date range" caption shared string (cell A2) ---
$ws.Range("A2").Value = "טווח תאריכים : 30/12/2018 - 31/12/2019"

# --- 3. Shift the existing rate table (rows 6:243) down by one row, then
#        fill in the new first data row (30/12/2018) ---
$data = $ws.Range("A6:B243").Value2
$ws.Range("A7:B244").Value2 = $data

$ws.Cells.Item(6, 1).Value2 = 43465
$ws.Cells.Item(6, 2).Value2 = 3.7480000000000002

# --- 4. Extend the table formatting down through the new trailing rows and
#        fill in the six additional dates (through 31/12/2019) ---
$ws.Range("A243:B243").Copy()
$ws.Range("A244:B250").PasteSpecial(-4122)

$ws.Cells.Item(245, 1).Value2 = 43822
$ws.Cells.Item(245, 2).Value2 = 3.472

$ws.Cells.Item(246, 1).Value2 = 43823
$ws.Cells.Item(246, 2).Value2 = 3.4660000000000002

$ws.Cells.Item(247, 1).Value2 = 43825
$ws.Cells.Item(247, 2).Value2 = 3.472

$ws.Cells.Item(248, 1).Value2 = 43826
$ws.Cells.Item(248, 2).Value2 = 3.468

$ws.Cells.Item(249, 1).Value2 = 43829
$ws.Cells.Item(249, 2).Value2 = 3.4630000000000001

$ws.Cells.Item(250, 1).Value2 = 43830
$ws.Cells.Item(250, 2).Value2 = 3.456

# --- 5. Reset the view: move the selection back to the top-left cell ---
$ws.Range("A1").Select()

Write-Output "edit complete"
